$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the longer "United Kingdom" / country labels
$ws.Columns.Item(1).ColumnWidth = 23

# --- New header row additions -------------------------------------------
$ws.Range("C1").Value = "Group"          # new shared string: Group

# --- Column C: a "Group" helper column -----------------------------------
# mostly literal "Other", with a handful of rows that mirror the country
# name back via a formula.
$ws.Range("C2").Formula = "=A2"
$ws.Range("C3").Value = "Other"          # new shared string: Other
$ws.Range("C4").Value = "Other"
$ws.Range("C5").Value = "Other"
$ws.Range("C6").Value = "Other"
$ws.Range("C7").Formula = "=A7"
$ws.Range("C8").Formula = "=A8"
$ws.Range("C9").Value = "Other"
$ws.Range("C10").Formula = "=A10"
$ws.Range("C11").Value = "Other"
$ws.Range("C12").Value = "Other"
$ws.Range("C13").Value = "Other"
$ws.Range("C14").Value = "Other"

# Fix the UK row label (drop the "(Pre-Brexit)" qualifier)
$ws.Range("A14").Value = "United Kingdom"   # new shared string: United Kingdom

$ws.Range("C15").Formula = "=A15"
$ws.Range("C16").Value = "Other"
$ws.Range("C17").Value = "Other"
$ws.Range("C18").Formula = "=A18"
$ws.Range("C19").Value = "Other"
$ws.Range("C20").Value = "Other"
$ws.Range("C21").Value = "Other"
$ws.Range("C22").Formula = "=A22"

# --- Column D: EU / NON-EU classification ---------------------------------
$ws.Range("D1").Value = "EU"             # new shared string: EU
$ws.Range("D2").Value = "NON-EU"         # new shared string: NON-EU
$ws.Range("D3").Value = "EU"
$ws.Range("D4").Value = "EU"
$ws.Range("D5").Value = "EU"
$ws.Range("D6").Value = "EU"
$ws.Range("D7").Value = "EU"
$ws.Range("D8").Value = "EU"
$ws.Range("D9").Value = "EU"
$ws.Range("D10").Value = "EU"
$ws.Range("D11").Value = "EU"
$ws.Range("D12").Value = "EU"
$ws.Range("D13").Value = "EU"
$ws.Range("D14").Value = "NON-EU"
$ws.Range("D15").Value = "EU"
$ws.Range("D16").Value = "EU"
$ws.Range("D17").Value = "EU"
$ws.Range("D18").Value = "EU"
$ws.Range("D19").Value = "EU"
$ws.Range("D20").Value = "EU"
$ws.Range("D21").Value = "EU"
$ws.Range("D22").Value = "EU"

$ws.Range("D14").Select() | Out-Null
